$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.117.49"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.908.38"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "567.01"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").Value = "143.74"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.905.04"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "0.432"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "32.50"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "3.391.03"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "62.056.14"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.56"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.906.23"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "430.76"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").Value = "13.04"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").Value = "0.652"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "78.71"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "11.96"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "25.65"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "5.38"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "48.87"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.113"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "8.14"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "40.85"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D44").Value = "2.724.19"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "133.27"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0337"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "344.27"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.000220"
$ws.Range("E50").Value = "  +13.26%  "
$ws.Range("E51").Value = "  -1.10%  "
